# Bugfixed evaluation and simulated rt_data for components.
# A new earliest data point (row) is inserted at the top of the forecast
# table, shifting all existing rows down by one, and the y_0_forecast /
# y_1_forecast (columns C and E) values are recomputed for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (the first data row),
# shifting all existing data rows down by one.
$ws.Rows.Item(2).Insert()

# The newly inserted row inherited generic formatting; restore the
# formatting (date style in column A, plain/general elsewhere) used by
# the rest of the data rows by copying the format from the row below
# (now row 3) onto it, then drop the now-unused column E placeholder
# since this row has no y_1_forecast figure (matching the other early
# rows that only started reporting column E from row 5 onward).
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E2").Clear()

$ws.Cells.Item(2, 1).Value2 = 39400
$ws.Cells.Item(2, 2).Value2 = 2007
$ws.Cells.Item(2, 3).Value2 = 11.13090654781821
$ws.Cells.Item(2, 4).Value2 = 2008
$ws.Cells.Item(3, 1).Value2 = 39765
$ws.Cells.Item(3, 2).Value2 = 2008
$ws.Cells.Item(3, 3).Value2 = 4.672550446571067
$ws.Cells.Item(3, 4).Value2 = 2009
$ws.Cells.Item(4, 1).Value2 = 40130
$ws.Cells.Item(4, 2).Value2 = 2009
$ws.Cells.Item(4, 3).Value2 = -14.45332333832743
$ws.Cells.Item(4, 4).Value2 = 2010
$ws.Cells.Item(5, 1).Value2 = 40494
$ws.Cells.Item(5, 2).Value2 = 2010
$ws.Cells.Item(5, 3).Value2 = 8.600536527919633
$ws.Cells.Item(5, 4).Value2 = 2011
$ws.Cells.Item(5, 5).Value2 = 9.002271992040312
$ws.Cells.Item(6, 1).Value2 = 40862
$ws.Cells.Item(6, 2).Value2 = 2011
$ws.Cells.Item(6, 3).Value2 = 10.25770250047622
$ws.Cells.Item(6, 4).Value2 = 2012
$ws.Cells.Item(6, 5).Value2 = 7.550992341868956
$ws.Cells.Item(7, 1).Value2 = 41228
$ws.Cells.Item(7, 2).Value2 = 2012
$ws.Cells.Item(7, 3).Value2 = 4.639893381363169
$ws.Cells.Item(7, 4).Value2 = 2013
$ws.Cells.Item(7, 5).Value2 = 5.799303245920884
$ws.Cells.Item(8, 1).Value2 = 41592
$ws.Cells.Item(8, 2).Value2 = 2013
$ws.Cells.Item(8, 3).Value2 = 0.3058963467304165
$ws.Cells.Item(8, 4).Value2 = 2014
$ws.Cells.Item(8, 5).Value2 = 2.638010271840896
$ws.Cells.Item(9, 1).Value2 = 41957
$ws.Cells.Item(9, 2).Value2 = 2014
$ws.Cells.Item(9, 3).Value2 = 4.068173739091874
$ws.Cells.Item(9, 4).Value2 = 2015
$ws.Cells.Item(9, 5).Value2 = 5.156937396706884
$ws.Cells.Item(10, 1).Value2 = 42321
$ws.Cells.Item(10, 2).Value2 = 2015
$ws.Cells.Item(10, 3).Value2 = 4.984288257750213
$ws.Cells.Item(10, 4).Value2 = 2016
$ws.Cells.Item(10, 5).Value2 = 3.4064284328156
$ws.Cells.Item(11, 1).Value2 = 42689
$ws.Cells.Item(11, 2).Value2 = 2016
$ws.Cells.Item(11, 3).Value2 = 1.878184267712912
$ws.Cells.Item(11, 4).Value2 = 2017
$ws.Cells.Item(11, 5).Value2 = 2.129835064860464
$ws.Cells.Item(12, 1).Value2 = 43053
$ws.Cells.Item(12, 2).Value2 = 2017
$ws.Cells.Item(12, 3).Value2 = 4.695933104194339
$ws.Cells.Item(12, 4).Value2 = 2018
$ws.Cells.Item(12, 5).Value2 = 5.022591279638045
$ws.Cells.Item(13, 1).Value2 = 43418
$ws.Cells.Item(13, 2).Value2 = 2018
$ws.Cells.Item(13, 3).Value2 = 4.892602738886098
$ws.Cells.Item(13, 4).Value2 = 2019
$ws.Cells.Item(13, 5).Value2 = 0.6513682883433347
$ws.Cells.Item(14, 1).Value2 = 43783
$ws.Cells.Item(14, 2).Value2 = 2019
$ws.Cells.Item(14, 3).Value2 = 0.8049382522247184
$ws.Cells.Item(14, 4).Value2 = 2020
$ws.Cells.Item(14, 5).Value2 = 3.036929265763488
$ws.Cells.Item(15, 1).Value2 = 44159
$ws.Cells.Item(15, 2).Value2 = 2020
$ws.Cells.Item(15, 3).Value2 = -8.784173899737169
$ws.Cells.Item(15, 4).Value2 = 2021
$ws.Cells.Item(15, 5).Value2 = 1.573231731123359
$ws.Cells.Item(16, 1).Value2 = 44525
$ws.Cells.Item(16, 2).Value2 = 2021
$ws.Cells.Item(16, 3).Value2 = 5.110501195359984
$ws.Cells.Item(16, 4).Value2 = 2022
$ws.Cells.Item(16, 5).Value2 = 1.88131082127776
$ws.Cells.Item(17, 1).Value2 = 44890
$ws.Cells.Item(17, 2).Value2 = 2022
$ws.Cells.Item(17, 3).Value2 = 5.120680133083599
$ws.Cells.Item(17, 4).Value2 = 2023
$ws.Cells.Item(17, 5).Value2 = 2.522545412785848
$ws.Cells.Item(18, 1).Value2 = 45254
$ws.Cells.Item(18, 2).Value2 = 2023
$ws.Cells.Item(18, 3).Value2 = -0.5532735011319234
$ws.Cells.Item(18, 4).Value2 = 2024
$ws.Cells.Item(18, 5).Value2 = -0.657715646732393
$ws.Cells.Item(19, 1).Value2 = 45618
$ws.Cells.Item(19, 2).Value2 = 2024
$ws.Cells.Item(19, 3).Value2 = -1.069674659641462
$ws.Cells.Item(19, 4).Value2 = 2025
$ws.Cells.Item(19, 5).Value2 = 0.5636794832278413
